# Week 15 simulations update
# Updates the "R" (Road) row (row 3) stats on both the OFF and DEF sheets.

$wb = $excel.ActiveWorkbook

# --- OFF sheet ---
$wsOff = $wb.Worksheets.Item("OFF")
$wsOff.Range("B3").Value = 205
$wsOff.Range("C3").Value = 149
$wsOff.Range("D3").Value = 48
$wsOff.Range("E3").Value = 34

# --- DEF sheet ---
$wsDef = $wb.Worksheets.Item("DEF")
$wsDef.Range("B3").Value = 356
$wsDef.Range("C3").Value = 263
$wsDef.Range("D3").Value = 90
$wsDef.Range("E3").Value = 44
